$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match-detail columns (F:V) between mis-ordered row pairs ---
# The "Indice"/date (A:E) stay put; only the match/odds data (F:V) needs swapping.
$swapPairs = @(
    @(42, 43),
    @(66, 67),
    @(70, 71)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("F$r1`:V$r1")
    $range2 = $ws.Range("F$r2`:V$r2")

    $tmp = $range1.Value2
    $range1.Value2 = $range2.Value2
    $range2.Value2 = $tmp
}

# --- Append three new match rows at the end (84:86) ---
# Copy formatting from the last existing row (83) so styles (bold/bordered
# index column, date-formatted E column) carry over correctly.
$ws.Range("A83:V83").Copy()
$ws.Range("A84:V86").PasteSpecial(-4122)  # xlPasteFormats

$newRows = @(
    @{
        A = 83; B = "portugal"; C = "liga-portugal"; D = "2023-2024"; E = 45234.6875
        F = "Chaves"; G = 0; H = "Benfica"; I = 2
        J = 10.21; K = "28/10/2023 19:13"
        L = 10.65; M = "04/11/2023 16:25"
        N = 6.93; O = "28/10/2023 19:13"
        P = 6.3; Q = "04/11/2023 16:25"
        R = 1.24; S = "28/10/2023 19:13"
        T = 1.28; U = "04/11/2023 16:25"
        V = "https://www.betexplorer.com/football/portugal/liga-portugal/chaves-benfica/ppCYXqTH/"
    },
    @{
        A = 84; B = "portugal"; C = "liga-portugal"; D = "2023-2024"; E = 45234.79166666666
        F = "Famalicao"; G = 3; H = "Gil Vicente"; I = 1
        J = 1.95; K = "29/10/2023 19:13"
        L = 2.16; M = "04/11/2023 18:53"
        N = 3.7; O = "29/10/2023 19:13"
        P = 3.47; Q = "04/11/2023 18:53"
        R = 4; S = "29/10/2023 19:13"
        T = 3.62; U = "04/11/2023 18:53"
        V = "https://www.betexplorer.com/football/portugal/liga-portugal/famalicao-gil-vicente/xK8QZ5b5/"
    },
    @{
        A = 85; B = "portugal"; C = "liga-portugal"; D = "2023-2024"; E = 45234.89583333334
        F = "Braga"; G = 6; H = "Portimonense"; I = 1
        J = 1.23; K = "30/10/2023 14:42"
        L = 1.21; M = "04/11/2023 21:26"
        N = 7.05; O = "30/10/2023 14:42"
        P = 7.62; Q = "04/11/2023 21:26"
        R = 12.14; S = "30/10/2023 14:42"
        T = 13.12; U = "04/11/2023 21:26"
        V = "https://www.betexplorer.com/football/portugal/liga-portugal/braga-portimonense/bT9Mznqa/"
    }
)

$rowNum = 84
foreach ($data in $newRows) {
    foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")) {
        $ws.Range("$col$rowNum").Value2 = $data[$col]
    }
    $rowNum = $rowNum + 1
}
